# Drools participant-assignment-rules.xlsx update
# Mirrors the xml_diff: Sheet1 rule-table text/content changes, a row-height
# tweak, clearing a stray duplicated cell, and moving the sheet's scroll /
# selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Imports block (rows 3-8): Spring expression classes shift up one slot
# because the now-unused "$acmObject: CheckParticipantListModel" / SpEL
# filter rows were removed above them.
$ws.Range("D3").Value = 'org.springframework.expression.EvaluationContext'
$ws.Range("D4").Value = 'org.springframework.expression.Expression'
$ws.Range("D5").Value = 'org.springframework.expression.ExpressionParser'
$ws.Range("D6").Value = 'org.springframework.expression.spel.standard.SpelExpressionParser'
$ws.Range("D7").Value = 'org.springframework.expression.spel.support.StandardEvaluationContext'
$ws.Range("D8").Value = 'com.armedia.acm.services.participants.model.CheckParticipantListModel'
# D9 (AcmParticipant import) is untouched.
$ws.Range("D10").Value = 'java.util.List'
$ws.Range("D11").Value = 'java.util.ArrayList'

# --- Functions block (row 12): rewritten Drools function source, shorter
# than before (the getAssignee/checkNoAccess helpers were dropped), so the
# wrapped row height shrinks from 409.5 to 360.
$funcText = @'
function Boolean evalBoolean(String expression, Object obj)
{
    ExpressionParser ep = new SpelExpressionParser();
    Expression exp = ep.parseExpression(expression);
    EvaluationContext ec = new StandardEvaluationContext();

    Boolean evaluated = exp.getValue(ec, obj, Boolean.class);
    return evaluated;
}
function void addErrorMessage(CheckParticipantListModel model, String error)
{
    if(error != null && !"".equals(error)){
         if(model.getErrorsList() == null)
         {
              model.setErrorsList(new ArrayList());
         }
          model.getErrorsList().add(error);
      }
}

'@
$ws.Range("D12").Value = $funcText
$ws.Rows.Item(12).RowHeight = 360

# --- Rule table header/condition/action rows (15-20).
$ws.Range("F16").Value = 'ACTION'

# Row 17: the model alias cell is renamed ($acmObject -> $model); the
# duplicated value that used to sit in F17 as well is removed entirely.
$ws.Range("C17").Value = '$model: CheckParticipantListModel'
$ws.Range("F17").ClearContents()

# Row 18: condition / action expression snippets, now operating on $model.
$ws.Range("C18").Value = 'objectType.equals("$param")'
$ws.Range("D18").Value = 'eval(evalBoolean("$param", $model))'
$ws.Range("E18").Value = 'eval(evalBoolean("$param", $model))'
$ws.Range("F18").Value = '$model.addErrorMessage("$param");'

# Row 19: column help text / captions.
$ws.Range("D19").Value = "Expression 1`n`nMust be a Spring expression that evaluates to true or false."
$ws.Range("E19").Value = "Expression 2`n`nMust be a Spring expression that evaluates to true or false. Check whether assignee is of not null"
$ws.Range("F19").Value = 'Add error to the list'

# Row 20: the actual rule row - new CASE_FILE assignee/no-access check.
$ws.Range("B20").Value = 'Case File - Check participants list for NoAccess & Owner'
$ws.Range("C20").Value = 'CASE_FILE'
$ws.Range("D20").Value = "participants != null && participants.containsKey('No Access') && participants.containsKey('assignee')"
$ws.Range("E20").Value = "participants['No Access'].contains(participants['assignee'][0])"
$ws.Range("F20").Value = 'Assignees cannot be on the no-access list.'

# E20 previously had the default (unstyled) format; match the bordered/filled
# style used by the rest of the row (same as copying D20's format onto it).
$ws.Range("D20").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- View state: scroll the sheet down and move the active selection to the
# bottom rule row.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D20").Select() | Out-Null
